$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 55.3746767090015
$ws.Range("E2").Value = 56.7933697165467
$ws.Range("L2").Value = 49.0485989036895
$ws.Range("B3").Value = 40.563158268296
$ws.Range("L3").Value = 36.9618773883661
